# Weekly update: a new Kiwi price record was reported for
# "Macroferia Regional de Talca" and needs to be inserted as a new data
# row right above the existing row 353, pushing the rest of the table
# (old rows 353-381) down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row above row 353; everything from 353 down shifts
# to 354.. (so old row 381 becomes row 382), matching the new
# dimension A1:T382.
$ws.Rows.Item(353).Insert()

# Fill the freshly inserted row 353 with the new observation. Columns
# A, B, C, E, F, G, H, I, J, K and R share the same market/product
# metadata as every other row in this block.
$ws.Range("A353").Value = 5
$ws.Range("B353").Value = "Macroferia Regional de Talca"
$ws.Range("C353").Value = "Maule"
$ws.Range("D353").Value = 44939
$ws.Range("E353").Value = 7
$ws.Range("F353").Value = "Fruta"
$ws.Range("G353").Value = 100101
$ws.Range("H353").Value = "Berries"
$ws.Range("I353").Value = 100101007
$ws.Range("J353").Value = "Kiwi"
$ws.Range("K353").Value = "Hayward"
$ws.Range("L353").Value = "Primera"
$ws.Range("M353").Value = 200
$ws.Range("N353").Value = 15000
$ws.Range("O353").Value = 15000
$ws.Range("P353").Value = 15000
$ws.Range("Q353").Value = "$/caja 15 kilos granel"
$ws.Range("R353").Value = "Provincia de Curicó"
$ws.Range("S353").Value = 1000
$ws.Range("T353").Value = 15
